$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---- Overview sheet ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $status
$ovw.Range("F2").Value = $status
$ovw.Columns.Item(5).ColumnWidth = 29.166666666666668
$ovw.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $status
$zhcn.Range("K2").Value = "2016-08-27 04:47:52"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $status
$dede.Range("K2").Value = "2016-08-27 04:47:59"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
